$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "[2023년 UPDATE] 머신러닝/딥러닝(PyTorch, TensorFlow) 최신 도커(docker) 업데이트 안내"
$ws.Range("E4").Value = "https://teddylee777.github.io/data_science/data-science-docker-2023"
$ws.Range("D6").Value = "[Python - 프로그래머스] 힙(Heap) > 디스크 컨트롤러"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/Python-%ED%94%84%EB%A1%9C%EA%B7%B8%EB%9E%98%EB%A8%B8%EC%8A%A4-%ED%9E%99Heap-%EB%94%94%EC%8A%A4%ED%81%AC-%EC%BB%A8%ED%8A%B8%EB%A1%A4%EB%9F%AC"
$ws.Range("D9").Value = "‘전문가를 행세하는 비전문가들’을 읽고"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/not-experts-faking-experts/#utm_source=rss&utm_medium=rss&utm_campaign=not-experts-faking-experts"
$ws.Range("D32").Value = "[파이썬 클린코드] Chapter1. 코드 포매팅과 도구"
$ws.Range("E32").Value = "https://dodonam.tistory.com/400"
$ws.Range("D51").Value = "[python] 제너레이터의 필요성(메모리 효율성)"
$ws.Range("E51").Value = "https://bskyvision.com/entry/python-%EC%A0%9C%EB%84%88%EB%A0%88%EC%9D%B4%ED%84%B0%EC%9D%98-%ED%95%84%EC%9A%94%EC%84%B1%EB%A9%94%EB%AA%A8%EB%A6%AC-%ED%9A%A8%EC%9C%A8%EC%84%B1"
